$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the table first so new columns pick up headers from the sheet when we set them below
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H3"))

# Rename header F1 from ">100 MW" to "100-500MW" and add the two new size-bucket headers
$ws.Range("F1").Value = "100-500MW"
$ws.Range("G1").Value = "500-1000MW"
$ws.Range("H1").Value = "> 1000MW"

# New data row 3 values
$ws.Range("G3").Value = 0.7
$ws.Range("H3").Value = 0.5

# New formulas row 2
$ws.Range("G2").Formula = "=1.9838*G3"
$ws.Range("H2").Formula = "=1.79032*H3"

# Match style of neighboring formula cells (s="1")
$ws.Range("F2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122)

# Move the active selection
$ws.Range("H4").Select() | Out-Null
